# Generate Report for Handback
# Updates the timestamp cells for the "4dec7601-e538-4226-a964-8c5ca7bab0a0"
# file's row (row 3) across the Overview, zh-cn and de-de sheets, reflecting
# a freshly (re-)generated handback report.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for row 3
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-16 16:41:55"

# zh-cn sheet: "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K) columns for row 3
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-16 16:41:50"
$wsZhCn.Range("K3").Value = "2016-08-16 16:42:17"

# de-de sheet: "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K) columns for row 3
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-16 16:41:55"
$wsDeDe.Range("K3").Value = "2016-08-16 16:42:24"
